$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 86 (shifts old rows 86-88 down to 88-90)
$ws.Rows.Item(86).Resize(2).Insert()

# Row 86: new weekly data
$ws.Range("A86").Value = 9
$ws.Range("B86").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C86").Value = "Metropolitana"
$ws.Range("D86").Value = 44746
$ws.Range("E86").Value = 13
$ws.Range("F86").Value = 100114002
$ws.Range("G86").Value = "Camote"
$ws.Range("H86").Value = "Sin especificar"
$ws.Range("I86").Value = "Primera"
$ws.Range("J86").Value = 610
$ws.Range("K86").Value = 11000
$ws.Range("L86").Value = 12000
$ws.Range("M86").Value = 11500
$ws.Range("N86").Value = "$/caja 18 kilos"
$ws.Range("O86").Value = "Perú"
$ws.Range("P86").Value = 639
$ws.Range("Q86").Value = 18
$ws.Range("R86").Value = "Hortaliza"

# Row 87: new weekly data
$ws.Range("A87").Value = 9
$ws.Range("B87").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C87").Value = "Metropolitana"
$ws.Range("D87").Value = 44746
$ws.Range("E87").Value = 13
$ws.Range("F87").Value = 100114002
$ws.Range("G87").Value = "Camote"
$ws.Range("H87").Value = "Sin especificar"
$ws.Range("I87").Value = "Primera"
$ws.Range("J87").Value = 1060
$ws.Range("K87").Value = 9000
$ws.Range("L87").Value = 10000
$ws.Range("M87").Value = 9500
$ws.Range("N87").Value = "$/malla 18 kilos"
$ws.Range("O87").Value = "Perú"
$ws.Range("P87").Value = 528
$ws.Range("Q87").Value = 18
$ws.Range("R87").Value = "Hortaliza"

Write-Host "done"
